# Slide 4 (sldId 263) - Title 1 (shape id 2): fix typo "effected" -> "affected"
# and join the two runs into one, keeping the 2nd run's formatting (dirty="0").
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# Put the full corrected sentence into the 2nd run (keeps its rPr, e.g. dirty="0"),
# then empty out the 1st run so only the corrected run remains.
$run2 = $tr.Runs(2)
$run2.Text = "Were overall prices affected?"
$run1 = $tr.Runs(1)
$run1.Text = ""
